# Update timestamps in the handback-status report, as produced by a
# fresh "Generate Report for Handback" run.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for the 4ac298d7... entry.
# This cell shares its text with de-de!H3 below, so both move together.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-09-01 15:04:31"

# zh-cn sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# for the 4ac298d7... entry.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-09-01 15:04:26"
$wsZhCn.Range("K3").Value = "2016-09-01 15:04:44"

# de-de sheet: "Correspond Handoff Datetime" (mirrors Overview!G3's text) and
# "Correspond Handback DateTime" for the 4ac298d7... entry.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-09-01 15:04:31"
$wsDeDe.Range("K3").Value = "2016-09-01 15:04:52"
